# Auto-generated Excel COM-interop script to apply scheduled-runner updates
# to the Hyperion_Profits workbook (Leve profitability calculations).
$wb = $excel.ActiveWorkbook

# ALC!row33: "Glazed and Confused" / "Clear Glass Lens"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 821.38464
$ws.Range("I33").Value = 875.5714
$ws.Range("J33").Value = 758.1667
$ws.Range("K33").Value = 875.5714
$ws.Range("L33").Value = 758.1667
$ws.Range("M33").Value = -646.5714
$ws.Range("N33").Value = -1216.1667

# ALC!row39: "Riches' Brew" / "Hi-Potion of Mind"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 216.6
$ws.Range("I39").Value = 99.75
$ws.Range("J39").Value = 350.14285
$ws.Range("K39").Value = 299.25
$ws.Range("L39").Value = 1050.42855
$ws.Range("M39").Value = -3.25
$ws.Range("N39").Value = -1642.42855

# ALC!row43: "Growing Is Knowing" / "Growth Formula Gamma"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 33334936
$ws.Range("I43").Value = 35715860
$ws.Range("J43").Value = 1992
$ws.Range("K43").Value = 35715860
$ws.Range("L43").Value = 1992
$ws.Range("M43").Value = -35715791
$ws.Range("N43").Value = -2130

# ALC!row44: "Alive and Unwell" / "Budding Oak Wand"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H44").Value = 1007535.5
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 1007535.5
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 1007535.5
$ws.Range("N44").Value = -1008459.5

# ALC!row51: "A Bile Business" / "Shark Oil"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7403.778
$ws.Range("I51").Value = 7240
$ws.Range("J51").Value = 7441
$ws.Range("K51").Value = 7240
$ws.Range("L51").Value = 7441
$ws.Range("M51").Value = -6756
$ws.Range("N51").Value = -8409

# ALC!row69: "Steeling the Knife, Steeling the Mind" / "Grade 1 Mind Dissolvent"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 45711440
$ws.Range("I69").Value = 2750000
$ws.Range("J69").Value = 50007584
$ws.Range("K69").Value = 8250000
$ws.Range("L69").Value = 150022752
$ws.Range("M69").Value = -8249126
$ws.Range("N69").Value = -150024500

# ALC!row72: "Surgical Substitution (L)" / "Grade 1 Mind Dissolvent"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 45711440
$ws.Range("I72").Value = 2750000
$ws.Range("J72").Value = 50007584
$ws.Range("K72").Value = 24750000
$ws.Range("L72").Value = 450068256
$ws.Range("M72").Value = -24745632
$ws.Range("N72").Value = -450076992

# ALC!row96: "Scroll Down" / "Grade 1 Reisui of Intelligence"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 591
$ws.Range("I96").Value = 570.7692
$ws.Range("J96").Value = 678.6667
$ws.Range("K96").Value = 1712.3076
$ws.Range("L96").Value = 2036.0001
$ws.Range("M96").Value = -339.3075999999999
$ws.Range("N96").Value = -4782.0001

# ALC!row138: "All-night Crafting" / "Cunning Craftsman's Tisane"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3826.6533
$ws.Range("I138").Value = 1507.8334
$ws.Range("J138").Value = 4268.3335
$ws.Range("K138").Value = 4523.5002
$ws.Range("L138").Value = 12805.0005
$ws.Range("M138").Value = 616.4997999999996
$ws.Range("N138").Value = -23085.0005

# ALC!row141: "Remedy for Reason" / "Grade 1 Gemdraught of Mind"
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 9413.324000000001
$ws.Range("I141").Value = 4408.2334
$ws.Range("J141").Value = 30863.715
$ws.Range("K141").Value = 13224.7002
$ws.Range("L141").Value = 92591.145
$ws.Range("M141").Value = -8044.700199999999
$ws.Range("N141").Value = -102951.145

# ARM!row32: "Ingot We Trust" / "Steel Ingot"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4856.74
$ws.Range("I32").Value = 2903.3794
$ws.Range("J32").Value = 10819.632
$ws.Range("K32").Value = 2903.3794
$ws.Range("L32").Value = 10819.632
$ws.Range("M32").Value = -2616.3794
$ws.Range("N32").Value = -11393.632

# ARM!row61: "Dealing with the Tough Stuff" / "Cobalt Ingot"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3885.5483
$ws.Range("I61").Value = 4583.1
$ws.Range("J61").Value = 3553.3809
$ws.Range("K61").Value = 4583.1
$ws.Range("L61").Value = 3553.3809
$ws.Range("M61").Value = -4371.1
$ws.Range("N61").Value = -3977.3809

# ARM!row136: "Metal with Mettle" / "Cobalt Tungsten Ingot"
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3885.5483
$ws.Range("I136").Value = 4583.1
$ws.Range("J136").Value = 3553.3809
$ws.Range("K136").Value = 13749.3
$ws.Range("L136").Value = 10660.1427
$ws.Range("M136").Value = -11199.3
$ws.Range("N136").Value = -15760.1427

# BSM!row38: "The Naked Blade" / "Steel Falchion"
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 32076.924
$ws.Range("I38").Value = 30000
$ws.Range("J38").Value = 39000
$ws.Range("K38").Value = 30000
$ws.Range("L38").Value = 39000
$ws.Range("M38").Value = -29584
$ws.Range("N38").Value = -39832

# BSM!row134: "Ruthenium Supremium" / "Ruthenium Ingot"
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3186.365
$ws.Range("I134").Value = 1120.9474
$ws.Range("J134").Value = 6325.8
$ws.Range("K134").Value = 3362.8422
$ws.Range("L134").Value = 18977.4
$ws.Range("M134").Value = -827.8422
$ws.Range("N134").Value = -24047.4

# CRP!row31: "Wall Not Found" / "Walnut Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16833.016
$ws.Range("I31").Value = 1849.3448
$ws.Range("J31").Value = 29248.057
$ws.Range("K31").Value = 1849.3448
$ws.Range("L31").Value = 29248.057
$ws.Range("M31").Value = -1554.3448
$ws.Range("N31").Value = -29838.057

# CRP!row34: "Armoires of the Rich and Famous" / "Walnut Lumber"
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 16833.016
$ws.Range("I34").Value = 1849.3448
$ws.Range("J34").Value = 29248.057
$ws.Range("K34").Value = 1849.3448
$ws.Range("L34").Value = 29248.057
$ws.Range("M34").Value = -1647.3448
$ws.Range("N34").Value = -29652.057

# CUL!row5: "What a Sap" / "Maple Syrup"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1323.375
$ws.Range("I5").Value = 826.82355
$ws.Range("J5").Value = 2529.2856
$ws.Range("K5").Value = 2480.47065
$ws.Range("L5").Value = 7587.8568
$ws.Range("M5").Value = -2368.47065
$ws.Range("N5").Value = -7811.8568

# CUL!row80: "Saucy for a Suitor" / "Hollandaise Sauce"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 10000
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 30000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -29064
$ws.Range("N80").ClearContents()

# CUL!row83: "Saved by the Sauce (L)" / "Hollandaise Sauce"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H83").Value = 10000
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 90000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -85320
$ws.Range("N83").ClearContents()

# CUL!row110: "His Dark Utensils" / "Spaghetti al Nero"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 15137.625
$ws.Range("I110").Value = 3700.3333
$ws.Range("J110").Value = 22000
$ws.Range("K110").Value = 11100.9999
$ws.Range("L110").Value = 66000
$ws.Range("M110").Value = -7010.999899999999
$ws.Range("N110").Value = -74180

# CUL!row111: "Soup for the Soldier" / "Broad Bean Soup"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H111").Value = 10106.714
$ws.Range("I111").Value = 1124.5
$ws.Range("J111").Value = 64000
$ws.Range("K111").Value = 3373.5
$ws.Range("L111").Value = 192000
$ws.Range("M111").Value = -306.5
$ws.Range("N111").Value = -198134

# CUL!row113: "Can't Eat Just One" / "Night Vinegar"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3569.96
$ws.Range("I113").Value = 10380
$ws.Range("J113").Value = 1867.45
$ws.Range("K113").Value = 31140
$ws.Range("L113").Value = 5602.35
$ws.Range("M113").Value = -28970
$ws.Range("N113").Value = -9942.35

# CUL!row119: "Super Dark Times" / "Risotto al Nero"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 3906.3
$ws.Range("I119").Value = 3906.3
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 11718.9
$ws.Range("L119").Value = 0
$ws.Range("M119").Value = -6880.900000000001

# CUL!row132: "More Mezcal" / "Cooking Mezcal"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2631.2368
$ws.Range("I132").Value = 845.2222
$ws.Range("J132").Value = 3185.5173
$ws.Range("K132").Value = 7606.999800000001
$ws.Range("L132").Value = 28669.6557
$ws.Range("M132").Value = -5076.999800000001
$ws.Range("N132").Value = -33729.6557

# CUL!row135: "Not-so-secret Ingredient" / "Royal Maple Syrup"
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1323.375
$ws.Range("I135").Value = 826.82355
$ws.Range("J135").Value = 2529.2856
$ws.Range("K135").Value = 7441.41195
$ws.Range("L135").Value = 22763.5704
$ws.Range("M135").Value = -4906.41195
$ws.Range("N135").Value = -27833.5704

# GSM!row70: "Sky Is the Limit" / "Mythrite Ingot"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9735.223
$ws.Range("I70").Value = 12120.25
$ws.Range("J70").Value = 4965.1665
$ws.Range("K70").Value = 12120.25
$ws.Range("L70").Value = 4965.1665
$ws.Range("M70").Value = -11850.25
$ws.Range("N70").Value = -5505.1665

# GSM!row73: "Hulls of Broken Dreams (L)" / "Mythrite Ingot"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9735.223
$ws.Range("I73").Value = 12120.25
$ws.Range("J73").Value = 4965.1665
$ws.Range("K73").Value = 12120.25
$ws.Range("L73").Value = 4965.1665
$ws.Range("M73").Value = -11184.25
$ws.Range("N73").Value = -6837.1665

# GSM!row80: "Needs More Prayerbell" / "Hardsilver Ingot"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 50202804
$ws.Range("I80").Value = 100002880
$ws.Range("J80").Value = 402730
$ws.Range("K80").Value = 100002880
$ws.Range("L80").Value = 402730
$ws.Range("M80").Value = -100001882
$ws.Range("N80").Value = -404726

# GSM!row83: "With a Noise That Reaches Heaven (L)" / "Hardsilver Ingot"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 50202804
$ws.Range("I83").Value = 100002880
$ws.Range("J83").Value = 402730
$ws.Range("K83").Value = 500014400
$ws.Range("L83").Value = 2013650
$ws.Range("M83").Value = -500009408
$ws.Range("N83").Value = -2023634

# GSM!row113: "Copious Crystal Cannons" / "Manasilver Nugget"
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2797.375
$ws.Range("I113").Value = 2797.375
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2797.375
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -627.375
$ws.Range("N113").ClearContents()

# LTW!row40: "Best Served Toad" / "Toad Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6781.16
$ws.Range("I40").Value = 5510.524
$ws.Range("J40").Value = 13452
$ws.Range("K40").Value = 5510.524
$ws.Range("L40").Value = 13452
$ws.Range("M40").Value = -5374.524
$ws.Range("N40").Value = -13724

# LTW!row55: "It's Not a Job, It's a Calling" / "Peiste Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1841.7059
$ws.Range("I55").Value = 2588.2856
$ws.Range("J55").Value = 1319.1
$ws.Range("K55").Value = 2588.2856
$ws.Range("L55").Value = 1319.1
$ws.Range("M55").Value = -2415.2856
$ws.Range("N55").Value = -1665.1

# LTW!row93: "Hide to Go Seek" / "Gagana Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1614.3125
$ws.Range("I93").Value = 1443
$ws.Range("J93").Value = 1834.5714
$ws.Range("K93").Value = 1443
$ws.Range("L93").Value = 1834.5714
$ws.Range("M93").Value = -195
$ws.Range("N93").Value = -4330.5714

# LTW!row100: "Tiger in the Sack" / "Tiger Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 52073.57
$ws.Range("I100").Value = 4677.25
$ws.Range("J100").Value = 1000000
$ws.Range("K100").Value = 4677.25
$ws.Range("L100").Value = 1000000
$ws.Range("M100").Value = -4136.25
$ws.Range("N100").Value = -1001082

# LTW!row132: "Tenets of Tanning" / "Silver Lobo Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11388.059
$ws.Range("I132").Value = 12938.23
$ws.Range("J132").Value = 6350
$ws.Range("K132").Value = 38814.69
$ws.Range("L132").Value = 19050
$ws.Range("M132").Value = -36284.69
$ws.Range("N132").Value = -24110

# LTW!row136: "Respect for Br'aax" / "Br'aax Leather"
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 32168.111
$ws.Range("I136").Value = 41502
$ws.Range("J136").Value = 7900
$ws.Range("K136").Value = 124506
$ws.Range("L136").Value = 23700
$ws.Range("M136").Value = -121956
$ws.Range("N136").Value = -28800

# WVR!row14: "Hat in Hand" / "Straw Hat"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 8089.6
$ws.Range("I14").Value = 8155.222
$ws.Range("J14").Value = 7499
$ws.Range("K14").Value = 8155.222
$ws.Range("L14").Value = 7499
$ws.Range("M14").Value = -7987.222
$ws.Range("N14").Value = -7835

# WVR!row37: "Bet You Anything" / "Velveteen Sarouel of Gathering"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 35307.668
$ws.Range("I37").Value = 35461.5
$ws.Range("J37").Value = 35000
$ws.Range("K37").Value = 35461.5
$ws.Range("L37").Value = 35000
$ws.Range("M37").Value = -35258.5
$ws.Range("N37").Value = -35406

# WVR!row45: "Private Concerns" / "Linen Trousers"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18491.6
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 18491.6
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 18491.6
$ws.Range("N45").Value = -19473.6
$ws.Range("M45").ClearContents()

# WVR!row53: "I'll Swap You" / "Ranger's Tunic"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("M53").ClearContents()

# WVR!row81: "Where the Dragonflies, the Net Catches" / "Crawler Silk"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2724.8
$ws.Range("I81").Value = 2938.7778
$ws.Range("J81").Value = 799
$ws.Range("K81").Value = 5877.5556
$ws.Range("L81").Value = 1598
$ws.Range("M81").Value = -4816.5556
$ws.Range("N81").Value = -3720

# WVR!row84: "To Kill a Dragon on Nameday (L)" / "Crawler Silk"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2724.8
$ws.Range("I84").Value = 2938.7778
$ws.Range("J84").Value = 799
$ws.Range("K84").Value = 29387.778
$ws.Range("L84").Value = 7990
$ws.Range("M84").Value = -24083.778
$ws.Range("N84").Value = -18598

# WVR!row96: "Skills on Display" / "Ruby Cotton Cloth"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 6490.364
$ws.Range("I96").Value = 5601.5
$ws.Range("J96").Value = 6998.2856
$ws.Range("K96").Value = 5601.5
$ws.Range("L96").Value = 6998.2856
$ws.Range("M96").Value = -4228.5
$ws.Range("N96").Value = -9744.285599999999

# WVR!row100: "Of Great Import" / "Kudzu Thread"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2342
$ws.Range("I100").Value = 1680
$ws.Range("J100").Value = 4990
$ws.Range("K100").Value = 3360
$ws.Range("L100").Value = 9980
$ws.Range("M100").Value = -2819
$ws.Range("N100").Value = -11062

# WVR!row125: "Color Coated" / "Almasty Serge Coat of Healing"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 99999
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 99999
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 99999
$ws.Range("N125").Value = -109839

# WVR!row132: "Comfy Cabins" / "Snow Cotton Cloth"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 31026.53
$ws.Range("I132").Value = 5746.16
$ws.Range("J132").Value = 101249.78
$ws.Range("K132").Value = 17238.48
$ws.Range("L132").Value = 303749.34
$ws.Range("M132").Value = -14708.48
$ws.Range("N132").Value = -308809.34

# WVR!row136: "Weaving the Envelope" / "Sarcenet Cloth"
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1869.2727
$ws.Range("I136").Value = 1783.5186
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 5350.5558
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -2800.5558
$ws.Range("N136").Value = -24600
